$wb = $excel.ActiveWorkbook

# Add the new "GS" worksheet after the last existing sheet (CMS)
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "GS"

# Header row (re-using shared strings already present in the workbook, including
# the trailing newline variants used by the other extract tabs)
$ws.Range("A1").Value = "Contact_ID`n"
$ws.Range("B1").Value = "Contact_Date`n"
$ws.Range("C1").Value = "Contact_Type_Code"
$ws.Range("D1").Value = "Contact_Type_Desc"
$ws.Range("E1").Value = "OM_Name`n"
$ws.Range("F1").Value = "OM_Key`n"
$ws.Range("G1").Value = "OM_Grade`n"
$ws.Range("H1").Value = "OM_Team_Key`n"
$ws.Range("I1").Value = "OM_Provider_Code`n"

# Row 2
$ws.Range("A2").Value = 999
$ws.Range("B2").Value = 43033
$ws.Range("C2").Value = "GS1"
$ws.Range("D2").Value = "Group supervision"
$ws.Range("E2").Value = "Billy Jones"
$ws.Range("F2").Value = 1234
$ws.Range("G2").Value = "PO"
$ws.Range("H2").Value = "T1"
$ws.Range("I2").Value = "ND01"

# Row 3
$ws.Range("A3").Value = 888
$ws.Range("B3").Value = 43033
$ws.Range("C3").Value = "GS1"
$ws.Range("D3").Value = "Group supervision"
$ws.Range("E3").Value = "Jane Jones"
$ws.Range("F3").Value = 5678
$ws.Range("G3").Value = "PO"
$ws.Range("H3").Value = "T1"
$ws.Range("I3").Value = "ND01"

# Row 4
$ws.Range("A4").Value = 777
$ws.Range("B4").Value = 43033
$ws.Range("C4").Value = "GS1"
$ws.Range("D4").Value = "Group supervision"
$ws.Range("E4").Value = "Thomas Boyle"
$ws.Range("F4").Value = 2468
$ws.Range("G4").Value = "PO"
$ws.Range("H4").Value = "T1"
$ws.Range("I4").Value = "ND01"

# Date format for column B (reuses the same date style used on the CMS sheet)
$ws.Range("B2:B4").NumberFormat = "m/d/yy"

# Undo the auto row-height bump that typing a newline into row 1 triggered
$ws.Rows.Item(1).AutoFit()

# Column widths - "best fit" widths matching each header/value's text, mirroring
# the source workbook's auto-fitted columns on this sheet
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 16.333333333333332
$ws.Columns.Item(5).ColumnWidth = 9.5
$ws.Columns.Item(6).ColumnWidth = 7.666666666666667
$ws.Columns.Item(7).ColumnWidth = 9.666666666666666
$ws.Columns.Item(8).ColumnWidth = 13.166666666666666
$ws.Columns.Item(9).ColumnWidth = 16.833333333333332

# Selection matching the final diff state for the GS sheet
$ws.Range("H16").Select()
